$d = $word.ActiveDocument

# 1) "O sistema exibe uma nova janela, ao clicar no botão novo, um formulário é mostrado referente ao cadastramento."
#    -> "O sistema exibe uma nova janela, um formulário referente a associação, dois campos serão disponibilizados,
#        para que o usuario ADM selecione o cliente e o exercício."
$r = $d.Content
$r.Find.Execute(", ao clicar no botão novo, um formulário é mostrado referente ao cadastramento", $false, $false, $false, $false, $false, $true, 1, $false, ", um formulário referente a associação, dois campos serão disponibilizados, para que o usuario ADM selecione o cliente e o exercício", 2) | Out-Null

# 2) "...informações da atribuição e clica no ícone de inclusão localizado no lado direito do registro."
#    -> "...informações da associação e clica no ícone de inclusão localizado no lado esquerdo do registro."
$r = $d.Content
$r.Find.Execute("da atribuição e clica no ícone de inclusão localizado no lado direito do registro", $false, $false, $false, $false, $false, $true, 1, $false, "da associação e clica no ícone de inclusão localizado no lado esquerdo do registro", 2) | Out-Null

# 3) "O sistema persiste as informações da atribuição." -> "...da associação."
$r = $d.Content
$r.Find.Execute("O sistema persiste as informações da atribuição.", $false, $false, $false, $false, $false, $true, 1, $false, "O sistema persiste as informações da associação.", 2) | Out-Null

# 5) "O sistema preenche o formulário com informações do pagamento selecionado."
#    -> "O sistema preenche o formulário com informações da associação selecionada."
$r = $d.Content
$r.Find.Execute("O sistema preenche o formulário com informações do pagamento selecionado.", $false, $false, $false, $false, $false, $true, 1, $false, "O sistema preenche o formulário com informações da associação selecionada.", 2) | Out-Null

# 7) "a atribuição é realizada" -> "a associação é realizada" (Pós-condições paragraph)
$r = $d.Content
$r.Find.Execute("a atribuição é realizada", $false, $false, $false, $false, $false, $true, 1, $false, "a associação é realizada", 2) | Out-Null

# 6) Move the _GoBack bookmark from the first empty paragraph (right after Fluxos Alternativo block)
#    to the second empty paragraph following the "Ao final..." (Pós-condições) paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute("Acadsystem.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$p1 = $r.Next(4, 1)
$p2 = $p1.Next(4, 1)
$target = $p2.Duplicate
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
